$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("日常分数+线上分数排名")
$ws.Range("F3:F12").Value = 1
